$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the dataset. It belongs
# chronologically "above" the existing row 209, so insert a fresh row at
# 209 (pushing the old 209..259 block down to 210..260) and fill it in
# with the new record's values. Columns that stay identical to the old
# row 209 (A,B,C,E,F,G,H,I,N,O,Q,R) are simply re-written too, since the
# inserted row starts out blank.
$ws.Rows.Item(209).Insert()

$ws.Range("A209").Value = 9
$ws.Range("B209").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C209").Value = "Metropolitana"
$ws.Range("D209").Value = 44641
$ws.Range("E209").Value = 13
$ws.Range("F209").Value = 100112043
$ws.Range("G209").Value = "Pepino ensalada"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 61
$ws.Range("K209").Value = 16000
$ws.Range("L209").Value = 18000
$ws.Range("M209").Value = 17016
$ws.Range("N209").Value = "`$/caja 60 unidades"
$ws.Range("O209").Value = "Región de Arica y Parinacota"
$ws.Range("P209").Value = 284
$ws.Range("Q209").Value = 60
$ws.Range("R209").Value = "Hortaliza"
